$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns with the latest
# crypto snapshot. Both columns are stored as plain text; numeric-looking
# Price values are prefixed with a leading apostrophe so Excel keeps them
# as text (preserving formatting such as "1.00" or "0.586") instead of
# silently converting them to numbers.

$ws.Range("D2").Value = "69.127.92"
$ws.Range("D3").Value = "3.369.81"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D5").Value = "'586.89"
$ws.Range("D6").Value = "'179.33"
$ws.Range("D10").Value = "'0.586"
$ws.Range("D11").Value = "'48.23"
$ws.Range("D12").Value = "'0.0000276"
$ws.Range("D13").Value = "'698.26"
$ws.Range("D14").Value = "3.920.35"
$ws.Range("D15").Value = "'8.52"
$ws.Range("D16").Value = "69.106.79"
$ws.Range("D17").Value = "3.380.48"
$ws.Range("D19").Value = "'17.62"
$ws.Range("D20").Value = "'11.28"
$ws.Range("D21").Value = "'0.900"
$ws.Range("D22").Value = "'5.48"
$ws.Range("D24").Value = "'101.52"
$ws.Range("D26").Value = "'2.72"
$ws.Range("D27").Value = "'9.60"
$ws.Range("D28").Value = "'33.48"
$ws.Range("D29").Value = "'8.62"
$ws.Range("D30").Value = "'7.04"
$ws.Range("D31").Value = "'11.15"
$ws.Range("D32").Value = "'555.62"
$ws.Range("D34").Value = "'3.54"
$ws.Range("D35").Value = "'58.17"
$ws.Range("D36").Value = "'1.00"
$ws.Range("D37").Value = "3.711.11"
$ws.Range("D38").Value = "'0.143"
$ws.Range("D39").Value = "'34.84"
$ws.Range("D40").Value = "'3.22"
$ws.Range("D41").Value = "'2.65"
$ws.Range("D42").Value = "'0.341"
$ws.Range("D44").Value = "'0.0417"
$ws.Range("D49").Value = "'1.34"
$ws.Range("D50").Value = "'131.94"
$ws.Range("D51").Value = "'2.63"

$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("E6").Value = "  +2.52%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +1.17%  "
$ws.Range("E9").Value = "  +4.84%  "
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("E11").Value = "  +6.28%  "
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("E13").Value = "  +5.79%  "
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("E20").Value = "  +3.29%  "
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("E22").Value = "  +2.46%  "
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("E24").Value = "  +3.44%  "
$ws.Range("E25").Value = "  +2.41%  "
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("E27").Value = "  +4.08%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  +2.69%  "
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("E34").Value = "  +9.92%  "
$ws.Range("E35").Value = "  +3.38%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("E38").Value = "  +9.62%  "
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("E40").Value = "  +3.81%  "
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("E44").Value = "  +3.06%  "
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("E50").Value = "  +3.73%  "
$ws.Range("E51").Value = "  -1.33%  "
